{"js": "// Highlight quantitative impact metrics (percentages, dollar amounts, etc.)\n// in bold + color (#2C3E50) across specific resume bullet points, matching\n// the \"hybrid bold + color highlighting for impact metrics\" commit.\n//\n// Strategy: for each target paragraph (identified by its un-highlighted\n// plain-text prefix), find the numeric/percentage/dollar tokens inside that\n// paragraph only (via Paragraph.search, case-sensitive, no wildcards so the\n// literal \"%\"/\"$\"/\"\u00b1\" characters are taken verbatim) and flip their font to\n// bold + the metric color. Word/Office.js automatically splits the run(s)\n// that text lives in, producing exactly the plain/bold/plain/... run\n// sequence described by the diff.\n\nconst HIGHLIGHT_COLOR = \"#2C3E50\";\n\n// Each entry: a unique prefix identifying the paragraph, and the ordered\n// list of literal substrings inside that paragraph to bold + color.\nconst TARGETS = [\n  {\n    prefix: \"\u2022 Discovered systematic race coding errors\",\n    tokens: [\"23%\", \"64%\"],\n  },\n  {\n    prefix: \"\u2022 Utilized advanced sampling methods\",\n    tokens: [\"\u00b14.2%\", \"\u00b12.1%\", \"71%\", \"87%\"],\n  },\n  {\n    prefix: \"\u2022 Trigonometric algorithm for boundary estimation\",\n    tokens: [\"73.5%\", \"$4.7M\"],\n  },\n  {\n    prefix: \"\u2022 Built real-time FEC analysis systems\",\n    tokens: [\"$2\"],\n  },\n  {\n    prefix: \"\u2022 Modernized legacy ETL processes\",\n    tokens: [\"57%\"],\n  },\n  {\n    prefix: \"\u2022 Algorithmic innovation: Pioneered trigonometric\",\n    tokens: [\"73.5%\"],\n  },\n  {\n    prefix: \"\u2022 $4.7M savings enabled\",\n    tokens: [\"$4.7M\"],\n  },\n  {\n    prefix: \"\u2022 178% accuracy improvement\",\n    tokens: [\"178%\"],\n  },\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (const target of TARGETS) {\n  const paragraph = paragraphs.items.find((p) => p.text.startsWith(target.prefix));\n  if (!paragraph) continue;\n\n  for (const token of target.tokens) {\n    const found = paragraph.search(token, { matchCase: true, matchWildcards: false });\n    found.load(\"items\");\n    await context.sync();\n\n    for (const hit of found.items) {\n      hit.font.bold = true;\n      hit.font.color = HIGHLIGHT_COLOR;\n    }\n    await context.sync();\n  }\n}\n", "ps1": "# Highlight quantitative impact metrics (percentages, dollar amounts, etc.)\n# in bold + color (#2C3E50) across specific resume bullet points, matching\n# the \"hybrid bold + color highlighting for impact metrics\" commit.\n#\n# Strategy: for each target paragraph (identified by a unique plain-text\n# substring), scope a Range to that paragraph and use Find.Execute (case-\n# sensitive, no wildcards) to locate each literal numeric/percentage/dollar\n# token in left-to-right order, advancing the search start past each hit.\n# Setting Font.Bold/Font.Color on the found (now-collapsed-to-the-hit) Range\n# causes Word to split the run(s) that text lives in, producing exactly the\n# plain/bold/plain/... run sequence described by the diff.\n\n$d = $word.ActiveDocument\n\n$HIGHLIGHT_COLOR = \"2C3E50\"\n\n# Each entry: a unique substring identifying the paragraph, and the ordered\n# list of literal substrings inside that paragraph to bold + color.\n$targets = @(\n    @{ Match = \"Discovered systematic race coding errors\"; Tokens = @(\"23%\", \"64%\") },\n    @{ Match = \"Utilized advanced sampling methods\"; Tokens = @(\"\u00b14.2%\", \"\u00b12.1%\", \"71%\", \"87%\") },\n    @{ Match = \"Trigonometric algorithm for boundary estimation\"; Tokens = @(\"73.5%\", \"$4.7M\") },\n    @{ Match = \"Built real-time FEC analysis systems\"; Tokens = @(\"$2\") },\n    @{ Match = \"Modernized legacy ETL processes\"; Tokens = @(\"57%\") },\n    @{ Match = \"Algorithmic innovation: Pioneered trigonometric\"; Tokens = @(\"73.5%\") },\n    @{ Match = \"$4.7M savings enabled\"; Tokens = @(\"$4.7M\") },\n    @{ Match = \"178% accuracy improvement\"; Tokens = @(\"178%\") }\n)\n\nforeach ($target in $targets) {\n    foreach ($p in $d.Paragraphs) {\n        $pText = $p.Range.Text\n        if ($pText.Contains($target.Match)) {\n            $pStart = $p.Range.Start\n            $pEnd = $p.Range.End\n            $searchFrom = $pStart\n\n            foreach ($tok in $target.Tokens) {\n                $fr = $d.Range($searchFrom, $pEnd)\n                $find = $fr.Find\n                $find.ClearFormatting()\n                $find.Text = $tok\n                $find.MatchCase = $true\n                $find.MatchWildcards = $false\n                $find.Forward = $true\n                $found = $find.Execute()\n                if ($found) {\n                    $fr.Font.Bold = 1\n                    $fr.Font.Color = $HIGHLIGHT_COLOR\n                    $searchFrom = $fr.End\n                }\n            }\n\n            break\n        }\n    }\n}\n"}
